$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header values in row 1 for new columns P and Q
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

for ($r = 2; $r -le 25; $r++) {
    # Swap values in I, K, M, O columns
    $ws.Range("I$r").Value = 2
    $ws.Range("K$r").Value = 1
    $ws.Range("M$r").Value = 2
    $ws.Range("O$r").Value = 1

    # Add new columns P and Q with value 2
    $ws.Range("P$r").Value = 2
    $ws.Range("Q$r").Value = 2
}
